$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.730.79'
$ws.Range("E2").Value = '  +2.26%  '

$ws.Range("D3").Value = '3.374.57'
$ws.Range("E3").Value = '  +0.95%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.48%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.79'
$ws.Range("E5").Value = '  +6.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.97'
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.601'
$ws.Range("E7").Value = '  +2.98%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.183'
$ws.Range("E9").Value = '  +2.80%  '

$ws.Range("E10").Value = '  +1.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.48'
$ws.Range("E11").Value = '  +2.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000278'
$ws.Range("E12").Value = '  +3.81%  '

$ws.Range("D13").Value = '3.918.70'
$ws.Range("E13").Value = '  +1.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '638.12'
$ws.Range("E14").Value = '  +8.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.62'
$ws.Range("E15").Value = '  +1.09%  '

$ws.Range("D16").Value = '67.768.97'
$ws.Range("E16").Value = '  +2.32%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.119'
$ws.Range("E17").Value = '  +1.25%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.375.30'
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.10'
$ws.Range("E19").Value = '  +1.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.14'
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.911'
$ws.Range("E21").Value = '  +1.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.97'
$ws.Range("E22").Value = '  -1.34%  '

$ws.Range("E23").Value = '  +1.73%  '

$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("E25").Value = '  +1.67%  '

$ws.Range("E26").Value = '  +5.92%  '

$ws.Range("E27").Value = '  +3.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.04'
$ws.Range("E28").Value = '  +7.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.71'
$ws.Range("E29").Value = '  +2.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.91'
$ws.Range("E30").Value = '  +4.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '613.19'
$ws.Range("E31").Value = '  +5.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.81'
$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("D33").Value = '4.055.30'
$ws.Range("E33").Value = '  +7.66%  '

$ws.Range("E34").Value = '  +1.74%  '

$ws.Range("E35").Value = '  +2.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.35'
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.80'
$ws.Range("E38").Value = '  +5.53%  '

$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '33.98'
$ws.Range("E39").Value = '  -1.10%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("E40").Value = '  +3.12%  '

$ws.Range("E41").Value = '  +2.84%  '

$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.41'
$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.344'
$ws.Range("E44").Value = '  +2.02%  '

$ws.Range("E45").Value = '  +2.48%  '

$ws.Range("E46").Value = '  +1.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.60'
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("E48").Value = '  +12.18%  '

$ws.Range("E49").Value = '  +0.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.26'
$ws.Range("E50").Value = '  +1.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.72'
$ws.Range("E51").Value = '  +4.63%  '
